# Update "Neg_Change" worksheet data (rows 2-8) and extend dimension to include new row 8.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")

$ws1.Range("A2").Value = "ONGC"
$ws1.Range("B2").Value = 273.89
$ws1.Range("C2").Value = 276.5
$ws1.Range("D2").Value = 268.01
$ws1.Range("E2").Value = 268.06
$ws1.Range("F2").Value = 24420043
$ws1.Range("G2").Value = 48267023
$ws1.Range("H2").Value = -0.4940636177209438
$ws1.Range("I2").Value = "ONGC"

$ws1.Range("A3").Value = "LODHA"
$ws1.Range("B3").Value = 947.7
$ws1.Range("C3").Value = 977
$ws1.Range("D3").Value = 936
$ws1.Range("E3").Value = 971.1
$ws1.Range("F3").Value = 2914884
$ws1.Range("G3").Value = 6125720
$ws1.Range("H3").Value = -0.5241565073166909
$ws1.Range("I3").Value = "LODHA"

$ws1.Range("A4").Value = "PNB"
$ws1.Range("B4").Value = 124.5
$ws1.Range("C4").Value = 126.03
$ws1.Range("D4").Value = 123.26
$ws1.Range("E4").Value = 125.1
$ws1.Range("F4").Value = 14202987
$ws1.Range("G4").Value = 32910164
$ws1.Range("H4").Value = -0.5684315945675628
$ws1.Range("I4").Value = "PNB"

$ws1.Range("A5").Value = "SIEMENS"
$ws1.Range("B5").Value = 3092.6
$ws1.Range("C5").Value = 3115
$ws1.Range("D5").Value = 3052.3
$ws1.Range("E5").Value = 3088
$ws1.Range("F5").Value = 242521
$ws1.Range("G5").Value = 487145
$ws1.Range("H5").Value = -0.5021584949039813
$ws1.Range("I5").Value = "SIEMENS"

$ws1.Range("A6").Value = "CANBK"
$ws1.Range("B6").Value = 150.7
$ws1.Range("C6").Value = 150.7
$ws1.Range("D6").Value = 146.62
$ws1.Range("E6").Value = 147.25
$ws1.Range("F6").Value = 43411151
$ws1.Range("G6").Value = 88706121
$ws1.Range("H6").Value = -0.5106183146031151
$ws1.Range("I6").Value = "CANBK"

$ws1.Range("A7").Value = "DALBHARAT"
$ws1.Range("B7").Value = 2060
$ws1.Range("C7").Value = 2075.1
$ws1.Range("D7").Value = 2032.9
$ws1.Range("E7").Value = 2069.9
$ws1.Range("F7").Value = 258368
$ws1.Range("G7").Value = 549292
$ws1.Range("H7").Value = -0.5296345113345907
$ws1.Range("I7").Value = "DALBHARAT"

$ws1.Range("A8").Value = "SRF"
$ws1.Range("B8").Value = 2818
$ws1.Range("C8").Value = 2840.4
$ws1.Range("D8").Value = 2788.4
$ws1.Range("E8").Value = 2806.1
$ws1.Range("F8").Value = 426581
$ws1.Range("G8").Value = 881502
$ws1.Range("H8").Value = -0.5160748359050802
$ws1.Range("I8").Value = "SRF"

# Update "Pos_Change" worksheet data (rows 2-21) and extend dimension to include new rows 18-21.
$ws2 = $wb.Worksheets.Item("Pos_Change")

$ws2.Range("A2").Value = "HINDUNILVR"
$ws2.Range("B2").Value = 2341
$ws2.Range("C2").Value = 2385
$ws2.Range("D2").Value = 2340
$ws2.Range("E2").Value = 2365
$ws2.Range("F2").Value = 2090462
$ws2.Range("G2").Value = 1478536
$ws2.Range("H2").Value = 0.4138729121238847
$ws2.Range("I2").Value = "HINDUNILVR"

$ws2.Range("A3").Value = "MAXHEALTH"
$ws2.Range("B3").Value = 945
$ws2.Range("C3").Value = 964
$ws2.Range("D3").Value = 940
$ws2.Range("E3").Value = 957.2
$ws2.Range("F3").Value = 5638365
$ws2.Range("G3").Value = 3538048
$ws2.Range("H3").Value = 0.5936372259505807
$ws2.Range("I3").Value = "MAXHEALTH"

$ws2.Range("A4").Value = "HDFCLIFE"
$ws2.Range("B4").Value = 721.25
$ws2.Range("C4").Value = 735.7
$ws2.Range("D4").Value = 720.65
$ws2.Range("E4").Value = 729.5
$ws2.Range("F4").Value = 3672455
$ws2.Range("G4").Value = 2511005
$ws2.Range("H4").Value = 0.4625438818321748
$ws2.Range("I4").Value = "HDFCLIFE"

$ws2.Range("A5").Value = "BAJAJFINSV"
$ws2.Range("B5").Value = 1928.5
$ws2.Range("C5").Value = 1956.6
$ws2.Range("D5").Value = 1927
$ws2.Range("E5").Value = 1953.7
$ws2.Range("F5").Value = 1342205
$ws2.Range("G5").Value = 884350
$ws2.Range("H5").Value = 0.5177305365522701
$ws2.Range("I5").Value = "BAJAJFINSV"

$ws2.Range("A6").Value = "BAJFINANCE"
$ws2.Range("B6").Value = 928.35
$ws2.Range("C6").Value = 938
$ws2.Range("D6").Value = 925.15
$ws2.Range("E6").Value = 930.85
$ws2.Range("F6").Value = 10151379
$ws2.Range("G6").Value = 6430445
$ws2.Range("H6").Value = 0.5786433131766153
$ws2.Range("I6").Value = "BAJFINANCE"

$ws2.Range("A7").Value = "GRASIM"
$ws2.Range("B7").Value = 2837.5
$ws2.Range("C7").Value = 2839.6
$ws2.Range("D7").Value = 2807.9
$ws2.Range("E7").Value = 2816.1
$ws2.Range("F7").Value = 1083570
$ws2.Range("G7").Value = 680590
$ws2.Range("H7").Value = 0.5921039098429304
$ws2.Range("I7").Value = "GRASIM"

$ws2.Range("A8").Value = "BOSCHLTD"
$ws2.Range("B8").Value = 35900
$ws2.Range("C8").Value = 36750
$ws2.Range("D8").Value = 35555
$ws2.Range("E8").Value = 36500
$ws2.Range("F8").Value = 24729
$ws2.Range("G8").Value = 16931
$ws2.Range("H8").Value = 0.4605752761207253
$ws2.Range("I8").Value = "BOSCHLTD"

$ws2.Range("A9").Value = "SOLARINDS"
$ws2.Range("B9").Value = 13228
$ws2.Range("C9").Value = 13564
$ws2.Range("D9").Value = 13010
$ws2.Range("E9").Value = 13443
$ws2.Range("F9").Value = 314201
$ws2.Range("G9").Value = 214858
$ws2.Range("H9").Value = 0.4623658416256318
$ws2.Range("I9").Value = "SOLARINDS"

$ws2.Range("A10").Value = "SHREECEM"
$ws2.Range("B10").Value = 27055
$ws2.Range("C10").Value = 27200
$ws2.Range("D10").Value = 26805
$ws2.Range("E10").Value = 26955
$ws2.Range("F10").Value = 42402
$ws2.Range("G10").Value = 28182
$ws2.Range("H10").Value = 0.5045773898232915
$ws2.Range("I10").Value = "SHREECEM"

$ws2.Range("A11").Value = "LTIM"
$ws2.Range("B11").Value = 5966
$ws2.Range("C11").Value = 5992.5
$ws2.Range("D11").Value = 5889
$ws2.Range("E11").Value = 5970
$ws2.Range("F11").Value = 260551
$ws2.Range("G11").Value = 182049
$ws2.Range("H11").Value = 0.4312135743673406
$ws2.Range("I11").Value = "LTIM"

$ws2.Range("A12").Value = "PFC"
$ws2.Range("B12").Value = 384
$ws2.Range("C12").Value = 384.45
$ws2.Range("D12").Value = 376
$ws2.Range("E12").Value = 378.5
$ws2.Range("F12").Value = 17579735
$ws2.Range("G12").Value = 11535966
$ws2.Range("H12").Value = 0.5239066238579413
$ws2.Range("I12").Value = "PFC"

$ws2.Range("A13").Value = "TORNTPOWER"
$ws2.Range("B13").Value = 1371.3
$ws2.Range("C13").Value = 1390
$ws2.Range("D13").Value = 1350.5
$ws2.Range("E13").Value = 1388
$ws2.Range("F13").Value = 559322
$ws2.Range("G13").Value = 366305
$ws2.Range("H13").Value = 0.5269297443387341
$ws2.Range("I13").Value = "TORNTPOWER"

$ws2.Range("A14").Value = "GODREJPROP"
$ws2.Range("B14").Value = 1568.1
$ws2.Range("C14").Value = 1588.6
$ws2.Range("D14").Value = 1550.1
$ws2.Range("E14").Value = 1576.1
$ws2.Range("F14").Value = 1924460
$ws2.Range("G14").Value = 1223856
$ws2.Range("H14").Value = 0.5724562366814396
$ws2.Range("I14").Value = "GODREJPROP"

$ws2.Range("A15").Value = "ASTRAL"
$ws2.Range("B15").Value = 1457
$ws2.Range("C15").Value = 1483.4
$ws2.Range("D15").Value = 1444.8
$ws2.Range("E15").Value = 1468
$ws2.Range("F15").Value = 371018
$ws2.Range("G15").Value = 231897
$ws2.Range("H15").Value = 0.5999258291396612
$ws2.Range("I15").Value = "ASTRAL"

$ws2.Range("A16").Value = "INDUSTOWER"
$ws2.Range("B16").Value = 440.95
$ws2.Range("C16").Value = 452
$ws2.Range("D16").Value = 435.45
$ws2.Range("E16").Value = 442.6
$ws2.Range("F16").Value = 11897285
$ws2.Range("G16").Value = 8487120
$ws2.Range("H16").Value = 0.4018047347038807
$ws2.Range("I16").Value = "INDUSTOWER"

$ws2.Range("A17").Value = "DABUR"
$ws2.Range("B17").Value = 512.15
$ws2.Range("C17").Value = 519.9
$ws2.Range("D17").Value = 501.35
$ws2.Range("E17").Value = 506.55
$ws2.Range("F17").Value = 2737132
$ws2.Range("G17").Value = 1887628
$ws2.Range("H17").Value = 0.4500378252494665
$ws2.Range("I17").Value = "DABUR"

$ws2.Range("A18").Value = "MUTHOOTFIN"
$ws2.Range("B18").Value = 4000
$ws2.Range("C18").Value = 4000
$ws2.Range("D18").Value = 3801.1
$ws2.Range("E18").Value = 3818
$ws2.Range("F18").Value = 1066990
$ws2.Range("G18").Value = 751686
$ws2.Range("H18").Value = 0.4194623818988248
$ws2.Range("I18").Value = "MUTHOOTFIN"

$ws2.Range("A19").Value = "HFCL"
$ws2.Range("B19").Value = 65.25
$ws2.Range("C19").Value = 69.75
$ws2.Range("D19").Value = 64.52
$ws2.Range("E19").Value = 67.90000000000001
$ws2.Range("F19").Value = 30274800
$ws2.Range("G19").Value = 20753501
$ws2.Range("H19").Value = 0.4587803763808333
$ws2.Range("I19").Value = "HFCL"

$ws2.Range("A20").Value = "BANDHANBNK"
$ws2.Range("B20").Value = 151.99
$ws2.Range("C20").Value = 155.49
$ws2.Range("D20").Value = 150.32
$ws2.Range("E20").Value = 154.45
$ws2.Range("F20").Value = 10411870
$ws2.Range("G20").Value = 7047526
$ws2.Range("H20").Value = 0.4773794378339292
$ws2.Range("I20").Value = "BANDHANBNK"

$ws2.Range("A21").Value = "MCX"
$ws2.Range("B21").Value = 2684
$ws2.Range("C21").Value = 2684
$ws2.Range("D21").Value = 2490
$ws2.Range("E21").Value = 2515
$ws2.Range("F21").Value = 7639727
$ws2.Range("G21").Value = 5308326
$ws2.Range("H21").Value = 0.4391970274621416
$ws2.Range("I21").Value = "MCX"
